$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ZS")

$ws.Range("B7").Value = 35000000.0
$ws.Range("C7").Value = 25000000.0
$ws.Range("D7").Value = 17916000.0
$ws.Range("E7").Value = 13326000.0
$ws.Range("F7").Value = 8541000.0
